$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, otherwise Excel auto-converts the string to a floating-point
# number (losing the exact text representation / trailing zeros).
$textCells = "D5","D6","D8","D10","D13","D16","D19","D21","D23","D25","D28","D29","D31","D33","D39","D40","D41","D44","D45","D49"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.188.02"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "2.354.36"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "540.29"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "135.99"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  +5.30%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "5.57"
$ws.Range("E10").Value = "  +5.25%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "23.83"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "2.773.22"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "58.197.11"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "0.0000134"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "2.352.00"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").Value = "332.65"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "62.83"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").Value = "8.50"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "1.75"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "171.97"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  +12.19%  "
$ws.Range("D33").Value = "18.45"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +6.96%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("D39").Value = "39.21"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "145.33"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").Value = "295.25"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").Value = "0.0946"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").Value = "19.26"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "17.55"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +0.52%  "
